$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-12-26 12:57:03"

# Column O ("timestamp") holds this value for every data row (2..547).
for ($r = 2; $r -le 547; $r++) {
    $ws.Cells.Item($r, 15).Value = $newTimestamp
}
